$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 7 for the new vendor "Carlisle Roofing Systems Inc."
# (vendor list is sorted alphabetically; Carlisle sits between
#  "Berridge Manufacturing Company" and "D & P Construction Co. Inc.")
$ws.Rows(7).Insert()

$ws.Range("A7").Value = "Carlisle Roofing Systems Inc."
$ws.Range("B7").Value = "x"
$ws.Range("G7").Value = "x"

# Update the "x" marks in column B (Envelope) to reflect the new state
$ws.Range("B4").Value = "x"
$ws.Range("B9").Value = ""
$ws.Range("B11").Value = ""
$ws.Range("B12").Value = "x"
$ws.Range("B15").Value = "x"
$ws.Range("B16").Value = "x"
$ws.Range("B17").Value = "x"
$ws.Range("B20").Value = ""

# Update the counter value used by monthly_bill_payment.py
$ws.Range("I1").Value = 33439

# Update the active selection to match the saved state
$ws.Range("B18").Select()
